# Automatische test-sync: 2025-07-22 12:41:50
#
# 1) "Logs" sheet: append row 14 with a new test-mail entry.
# 2) "Logs" sheet: extend the conditional-formatting ranges (D/G/H/I/J)
#    so row 14 is covered too.
# 3) "Dashboard" sheet: swap the category labels in rows 4 and 5, and add
#    a new row 7 for "Intern verzoek / Actie voor medewerker".
# 4) Chart1 (on Dashboard): extend the category/value series references
#    from row 6 to row 7.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet - add row 14
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Graag contact opnemen met de klant hierover."
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Testmail #14: Graag contact opnemen met de klant hierover."
$logs.Range("D14").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E14").Value = "Beste klant,`nBedankt voor uw bericht. Kunt u meer informatie geven over waarover u graag contact wilt opnemen? Bijvoorbeeld over welke specifieke kwestie of vraag het gaat? Met deze aanvullende details kunnen we u beter van dienst zijn.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Range("F14").Value = "2025-07-22 12:41:45"
$logs.Range("G14").Value = "Ja"
$logs.Range("H14").Value = "Nee"
$logs.Range("I14").Value = "Ja"
$logs.Range("J14").Value = "Ja"

# ---------------------------------------------------------------------
# 2) Logs sheet - grow the conditional formatting sqref from row 13 to
#    row 14 for columns D, G, H, I, J (rules/dxfIds stay unchanged).
# ---------------------------------------------------------------------
function Expand-ConditionalFormatRange($sheet, $col) {
    $oldRange = $sheet.Range($col + "2:" + $col + "13")
    $newRange = $sheet.Range($col + "2:" + $col + "14")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

Expand-ConditionalFormatRange $logs "D"
Expand-ConditionalFormatRange $logs "G"
Expand-ConditionalFormatRange $logs "H"
Expand-ConditionalFormatRange $logs "I"
Expand-ConditionalFormatRange $logs "J"

# ---------------------------------------------------------------------
# 3) Dashboard sheet - swap rows 4/5 labels, add new row 7
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("A5").Value = "Openingstijden / Locatie"

$dash.Range("A7").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B7").Value = 1

# ---------------------------------------------------------------------
# 4) Chart1 - extend category/value series references to row 7
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$7,Dashboard!`$B`$2:`$B`$7,1)"
